$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- View / window state -------------------------------------------------
# Reflect the updated viewport captured in the saved sheet view: new zoom
# level and scroll position, with the final selection left on I72 (the
# last cell the author edited in this session).
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1

# --- Data entry ------------------------------------------------------------
# Fill in the previously-empty "tamanho bytes" (F:I) and "Tempo total" (S:V)
# columns for each of the four result tables (rows 41-45, 50-54, 59-63) plus
# the last table's F:I columns (rows 68-72).
$ws.Range("F41").Value = 1590
$ws.Range("G41").Value = 6233
$ws.Range("H41").Value = 10496
$ws.Range("I41").Value = 17114
$ws.Range("S41").Value = 38561
$ws.Range("T41").Value = 37782
$ws.Range("U41").Value = 42085
$ws.Range("V41").Value = 46470
$ws.Range("F42").Value = 966
$ws.Range("G42").Value = 1672
$ws.Range("H42").Value = 10408
$ws.Range("I42").Value = 17107
$ws.Range("S42").Value = 39932
$ws.Range("T42").Value = 36408
$ws.Range("U42").Value = 42673
$ws.Range("V42").Value = 41116
$ws.Range("F43").Value = 1041
$ws.Range("G43").Value = 1097
$ws.Range("H43").Value = 10399
$ws.Range("I43").Value = 17196
$ws.Range("S43").Value = 35911
$ws.Range("T43").Value = 38903
$ws.Range("U43").Value = 38478
$ws.Range("V43").Value = 42328
$ws.Range("F44").Value = 977
$ws.Range("G44").Value = 407
$ws.Range("H44").Value = 10305
$ws.Range("I44").Value = 17302
$ws.Range("S44").Value = 34951
$ws.Range("T44").Value = 40298
$ws.Range("U44").Value = 37710
$ws.Range("V44").Value = 42204
$ws.Range("F45").Value = 943
$ws.Range("G45").Value = 384
$ws.Range("H45").Value = 10320
$ws.Range("I45").Value = 17628
$ws.Range("S45").Value = 38333
$ws.Range("T45").Value = 38256
$ws.Range("U45").Value = 45578
$ws.Range("V45").Value = 45541
$ws.Range("F50").Value = 18700
$ws.Range("G50").Value = 26053
$ws.Range("H50").Value = 33352
$ws.Range("I50").Value = 44839
$ws.Range("S50").Value = 4142
$ws.Range("T50").Value = 4731
$ws.Range("U50").Value = 5811
$ws.Range("V50").Value = 7526
$ws.Range("F51").Value = 2740
$ws.Range("G51").Value = 4434
$ws.Range("H51").Value = 32693
$ws.Range("I51").Value = 44177
$ws.Range("S51").Value = 4119
$ws.Range("T51").Value = 4870
$ws.Range("U51").Value = 5879
$ws.Range("V51").Value = 7580
$ws.Range("F52").Value = 2502
$ws.Range("G52").Value = 3433
$ws.Range("H52").Value = 32007
$ws.Range("I52").Value = 42667
$ws.Range("S52").Value = 4135
$ws.Range("T52").Value = 4790
$ws.Range("U52").Value = 5743
$ws.Range("V52").Value = 7028
$ws.Range("F53").Value = 2245
$ws.Range("G53").Value = 2090
$ws.Range("H53").Value = 33161
$ws.Range("I53").Value = 43335
$ws.Range("S53").Value = 4152
$ws.Range("T53").Value = 4813
$ws.Range("U53").Value = 5922
$ws.Range("V53").Value = 7362
$ws.Range("F54").Value = 2003
$ws.Range("G54").Value = 1391
$ws.Range("H54").Value = 35912
$ws.Range("I54").Value = 43456
$ws.Range("S54").Value = 4176
$ws.Range("T54").Value = 5149
$ws.Range("U54").Value = 5917
$ws.Range("V54").Value = 7372
$ws.Range("F59").Value = 9312
$ws.Range("G59").Value = 6233
$ws.Range("H59").Value = 13119
$ws.Range("I59").Value = 16821
$ws.Range("S59").Value = 800424
$ws.Range("T59").Value = 805175
$ws.Range("U59").Value = 883330
$ws.Range("V59").Value = 850585
$ws.Range("F60").Value = 553
$ws.Range("G60").Value = 1040
$ws.Range("H60").Value = 12835
$ws.Range("I60").Value = 16120
$ws.Range("S60").Value = 643701
$ws.Range("T60").Value = 633588
$ws.Range("U60").Value = 844258
$ws.Range("V60").Value = 841710
$ws.Range("F61").Value = 227
$ws.Range("G61").Value = 654
$ws.Range("H61").Value = 12421
$ws.Range("I61").Value = 15047
$ws.Range("S61").Value = 657475
$ws.Range("T61").Value = 663476
$ws.Range("U61").Value = 825944
$ws.Range("V61").Value = 834305
$ws.Range("F62").Value = 260
$ws.Range("G62").Value = 501
$ws.Range("H62").Value = 13349
$ws.Range("I62").Value = 15280
$ws.Range("S62").Value = 651661
$ws.Range("T62").Value = 648536
$ws.Range("U62").Value = 816126
$ws.Range("V62").Value = 830899
$ws.Range("F63").Value = 160
$ws.Range("G63").Value = 213
$ws.Range("H63").Value = 14499
$ws.Range("I63").Value = 15204
$ws.Range("S63").Value = 674228
$ws.Range("T63").Value = 671995
$ws.Range("U63").Value = 904180
$ws.Range("V63").Value = 841676
$ws.Range("F68").Value = 2915
$ws.Range("G68").Value = 3456
$ws.Range("H68").Value = 4179
$ws.Range("I68").Value = 5522
$ws.Range("F69").Value = 378
$ws.Range("G69").Value = 647
$ws.Range("H69").Value = 4131
$ws.Range("I69").Value = 5640
$ws.Range("F70").Value = 302
$ws.Range("G70").Value = 545
$ws.Range("H70").Value = 4081
$ws.Range("I70").Value = 5293
$ws.Range("F71").Value = 196
$ws.Range("G71").Value = 326
$ws.Range("H71").Value = 4057
$ws.Range("I71").Value = 5719
$ws.Range("F72").Value = 136
$ws.Range("G72").Value = 198
$ws.Range("H72").Value = 4164

# I72 additionally picks up an explicit font application (matching the
# author's last edit in that cell), which is why it ends up on its own
# cell style distinct from the rest of the column.
$ws.Range("I72").Value = 5440
$ws.Range("I72").Font.Name = "Calibri"

$ws.Range("I72").Select()
